# Auto-generated Excel COM-interop edit script
# Applies updated market-board price data to the Leve profit tables
# across all 8 job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 116.07143
$ws.Range("I5").Value = 110.416664
$ws.Range("J5").Value = 150
$ws.Range("K5").Value = 110.416664
$ws.Range("L5").Value = 150
$ws.Range("M5").Value = 4.583336000000003
$ws.Range("N5").Value = -380

$ws.Range("H32").Value = 843.381
$ws.Range("I32").Value = 547.5
$ws.Range("J32").Value = 913
$ws.Range("K32").Value = 547.5
$ws.Range("L32").Value = 913
$ws.Range("M32").Value = -221.5
$ws.Range("N32").Value = -1565

$ws.Range("H45").Value = 5466.6665
$ws.Range("I45").Value = 5400
$ws.Range("J45").Value = 5500
$ws.Range("K45").Value = 16200
$ws.Range("L45").Value = 16500
$ws.Range("M45").Value = -16008
$ws.Range("N45").Value = -16884

$ws.Range("H103").Value = 717.3077
$ws.Range("I103").Value = 717.3077
$ws.Range("K103").Value = 2151.9231
$ws.Range("M103").Value = -1565.9231

$ws.Range("H129").Value = 950608.4399999999
$ws.Range("I129").Value = 337
$ws.Range("J129").Value = 1123385
$ws.Range("K129").Value = 1011
$ws.Range("L129").Value = 3370155
$ws.Range("M129").Value = 3989
$ws.Range("N129").Value = -3380155

$ws.Range("H141").Value = 5726.25
$ws.Range("I141").Value = 6666.6665
$ws.Range("J141").Value = 2905
$ws.Range("K141").Value = 19999.9995
$ws.Range("L141").Value = 8715
$ws.Range("M141").Value = -14819.9995
$ws.Range("N141").Value = -19075

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H37").Value = 7160.2856
$ws.Range("I37").Value = 2000
$ws.Range("J37").Value = 8020.3335
$ws.Range("K37").Value = 2000
$ws.Range("L37").Value = 8020.3335
$ws.Range("M37").Value = -1727
$ws.Range("N37").Value = -8566.333500000001

$ws.Range("H61").Value = 1360
$ws.Range("I61").Value = 1312.8292
$ws.Range("K61").Value = 1312.8292
$ws.Range("M61").Value = -1100.8292

$ws.Range("H136").Value = 1360
$ws.Range("I136").Value = 1312.8292
$ws.Range("K136").Value = 3938.487599999999
$ws.Range("M136").Value = -1388.487599999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 260
$ws.Range("I22").Value = 260
$ws.Range("K22").Value = 260
$ws.Range("M22").Value = -87

$ws.Range("H86").Value = 27780612
$ws.Range("I86").Value = 31252442
$ws.Range("J86").Value = 5975
$ws.Range("K86").Value = 31252442
$ws.Range("L86").Value = 5975
$ws.Range("M86").Value = -31251319
$ws.Range("N86").Value = -8221

$ws.Range("H89").Value = 27780612
$ws.Range("I89").Value = 31252442
$ws.Range("J89").Value = 5975
$ws.Range("K89").Value = 156262210
$ws.Range("L89").Value = 29875
$ws.Range("M89").Value = -156256594
$ws.Range("N89").Value = -41107

$ws.Range("H105").Value = 8219.091
$ws.Range("I105").Value = 6888.75
$ws.Range("J105").Value = 11766.667
$ws.Range("K105").Value = 6888.75
$ws.Range("L105").Value = 11766.667
$ws.Range("M105").Value = -5141.75
$ws.Range("N105").Value = -15260.667

$ws.Range("H107").Value = 14066.3
$ws.Range("I107").Value = 857.875
$ws.Range("K107").Value = 857.875
$ws.Range("M107").Value = 1062.125

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 514.2857
$ws.Range("I7").Value = 558.3333
$ws.Range("J7").Value = 250
$ws.Range("K7").Value = 558.3333
$ws.Range("L7").Value = 250
$ws.Range("M7").Value = -445.3333
$ws.Range("N7").Value = -476

$ws.Range("H33").Value = 1281.909
$ws.Range("I33").Value = 460.1111
$ws.Range("K33").Value = 460.1111
$ws.Range("M33").Value = -81.11110000000002

$ws.Range("H107").Value = 483.9375
$ws.Range("I107").Value = 605.4
$ws.Range("J107").Value = 428.72726
$ws.Range("K107").Value = 605.4
$ws.Range("L107").Value = 428.72726
$ws.Range("M107").Value = 1314.6
$ws.Range("N107").Value = -4268.72726

$ws.Range("H122").Value = 1279.5
$ws.Range("I122").Value = 964
$ws.Range("J122").Value = 1595
$ws.Range("K122").Value = 2892
$ws.Range("L122").Value = 4785
$ws.Range("M122").Value = -442
$ws.Range("N122").Value = -9685

$ws.Range("H132").Value = 1356.1666
$ws.Range("I132").Value = 852.13794
$ws.Range("J132").Value = 3444.2856
$ws.Range("K132").Value = 2556.41382
$ws.Range("L132").Value = 10332.8568
$ws.Range("M132").Value = -26.41381999999976
$ws.Range("N132").Value = -15392.8568

$ws.Range("H134").Value = 17858392
$ws.Range("I134").Value = 1343.8667
$ws.Range("J134").Value = 38462676
$ws.Range("K134").Value = 4031.6001
$ws.Range("L134").Value = 115388028
$ws.Range("M134").Value = -1496.6001
$ws.Range("N134").Value = -115393098

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 2915.5293
$ws.Range("I137").Value = 2804.2856
$ws.Range("J137").Value = 2993.4
$ws.Range("K137").Value = 8412.856800000001
$ws.Range("L137").Value = 8980.200000000001
$ws.Range("M137").Value = -3312.856800000001
$ws.Range("N137").Value = -19180.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 9096527
$ws.Range("I11").Value = 18005920
$ws.Range("J11").Value = 1672033.4
$ws.Range("K11").Value = 18005920
$ws.Range("L11").Value = 1672033.4
$ws.Range("M11").Value = -18005781
$ws.Range("N11").Value = -1672311.4

$ws.Range("H12").Value = 70003
$ws.Range("I12").Value = 30000
$ws.Range("J12").Value = 83337.336
$ws.Range("K12").Value = 30000
$ws.Range("L12").Value = 83337.336
$ws.Range("M12").Value = -29860
$ws.Range("N12").Value = -83617.336

$ws.Range("H21").Value = 1000
$ws.Range("I21").Value = 1000
$ws.Range("K21").Value = 1000
$ws.Range("M21").Value = -827

$ws.Range("H30").Value = 1000
$ws.Range("I30").Value = 1000
$ws.Range("K30").Value = 1000
$ws.Range("M30").Value = -895

$ws.Range("H70").Value = 4721.4443
$ws.Range("I70").Value = 4844.5386
$ws.Range("K70").Value = 4844.5386
$ws.Range("M70").Value = -4574.5386

$ws.Range("H73").Value = 4721.4443
$ws.Range("I73").Value = 4844.5386
$ws.Range("K73").Value = 4844.5386
$ws.Range("M73").Value = -3908.5386

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 25000
$ws.Range("J20").Value = 25000
$ws.Range("L20").Value = 25000
$ws.Range("N20").Value = -25452

$ws.Range("H36").Value = 54715
$ws.Range("J36").Value = 54715
$ws.Range("L36").Value = 54715
$ws.Range("N36").Value = -55839

$ws.Range("H100").Value = 13890776
$ws.Range("I100").Value = 27780128
$ws.Range("J100").Value = 1425
$ws.Range("K100").Value = 27780128
$ws.Range("L100").Value = 1425
$ws.Range("M100").Value = -27779587
$ws.Range("N100").Value = -2507

$ws.Range("H136").Value = 3225.2593
$ws.Range("I136").Value = 1967.2222
$ws.Range("J136").Value = 5741.3335
$ws.Range("K136").Value = 5901.6666
$ws.Range("L136").Value = 17224.0005
$ws.Range("M136").Value = -3351.6666
$ws.Range("N136").Value = -22324.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1326
$ws.Range("I122").Value = 1326
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 3978
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -1528
$ws.Range("N122").Value = $null
